$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Productdata")

$ws.Range("C7").Value = 16001
$ws.Range("C8").Value = 8001
$ws.Range("C9").Value = 32001
$ws.Range("C10").Value = 2801
$ws.Range("C11").Value = 16001
$ws.Range("C12").Value = 8001
$ws.Range("C13").Value = 32001
$ws.Range("C14").Value = 34801
